$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price strings (e.g. "63.991.41", "2.20").
# Temporarily force text format on the data range so Excel keeps them as
# literal text instead of re-parsing them as numbers and dropping digits.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "63.991.41"
$ws.Range("E2").Value = "  +0.92%  "
# Row 3
$ws.Range("D3").Value = "3.139.65"
$ws.Range("E3").Value = "  +1.47%  "
# Row 4
$ws.Range("E4").Value = "  +0.09%  "
# Row 5
$ws.Range("D5").Value = "590.17"
$ws.Range("E5").Value = "  +1.47%  "
# Row 6
$ws.Range("D6").Value = "146.11"
$ws.Range("E6").Value = "  +1.11%  "
# Row 8
$ws.Range("D8").Value = "3.128.95"
$ws.Range("E8").Value = "  +1.30%  "
# Row 9
$ws.Range("E9").Value = "  +0.55%  "
# Row 10
$ws.Range("E10").Value = "  +1.85%  "
# Row 11
$ws.Range("D11").Value = "5.92"
$ws.Range("E11").Value = "  +5.35%  "
# Row 12
$ws.Range("D12").Value = "0.456"
$ws.Range("E12").Value = "  +0.23%  "
# Row 13
$ws.Range("D13").Value = "0.0000246"
$ws.Range("E13").Value = "  +0.49%  "
# Row 14
$ws.Range("D14").Value = "37.15"
$ws.Range("E14").Value = "  -0.77%  "
# Row 15
$ws.Range("D15").Value = "3.659.30"
$ws.Range("E15").Value = "  +1.48%  "
# Row 16
$ws.Range("E16").Value = "  -0.16%  "
# Row 17
$ws.Range("D17").Value = "7.26"
$ws.Range("E17").Value = "  +2.51%  "
# Row 18
$ws.Range("D18").Value = "63.794.11"
$ws.Range("E18").Value = "  +0.80%  "
# Row 19
$ws.Range("D19").Value = "3.136.93"
$ws.Range("E19").Value = "  +1.49%  "
# Row 20
$ws.Range("D20").Value = "467.56"
$ws.Range("E20").Value = "  +1.56%  "
# Row 21
$ws.Range("D21").Value = "14.40"
$ws.Range("E21").Value = "  +1.31%  "
# Row 22
$ws.Range("D22").Value = "0.730"
$ws.Range("E22").Value = "  +0.84%  "
# Row 23
$ws.Range("D23").Value = "7.54"
$ws.Range("E23").Value = "  +1.26%  "
# Row 24
$ws.Range("D24").Value = "2.38"
$ws.Range("E24").Value = "  +12.31%  "
# Row 25
$ws.Range("D25").Value = "13.09"
$ws.Range("E25").Value = "  +1.08%  "
# Row 26
$ws.Range("E26").Value = "  -0.33%  "
# Row 27
$ws.Range("E27").Value = "  +0.05%  "
# Row 28
$ws.Range("D28").Value = "9.88"
$ws.Range("E28").Value = "  +10.19%  "
# Row 29
$ws.Range("E29").Value = "  +1.65%  "
# Row 30
$ws.Range("D30").Value = "7.30"
$ws.Range("E30").Value = "  +6.90%  "
# Row 31
$ws.Range("E31").Value = "  +0.13%  "
# Row 32
$ws.Range("D32").Value = "2.20"
$ws.Range("E32").Value = "  +0.11%  "
# Row 33
$ws.Range("D33").Value = "0.115"
$ws.Range("E33").Value = "  +4.58%  "
# Row 34
$ws.Range("D34").Value = "27.63"
$ws.Range("E34").Value = "  +3.86%  "
# Row 35
$ws.Range("D35").Value = "0.0₃0854"
$ws.Range("E35").Value = "  +0.77%  "
# Row 36
$ws.Range("E36").Value = "  +2.88%  "
# Row 37
$ws.Range("D37").Value = "6.16"
$ws.Range("E37").Value = "  +2.97%  "
# Row 38
$ws.Range("D38").Value = "2.28"
$ws.Range("E38").Value = "  -0.70%  "
# Row 39
$ws.Range("E39").Value = "  -2.67%  "
# Row 40
$ws.Range("B40").Value = "Cosmos"
$ws.Range("C40").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D40").Value = "9.36"
$ws.Range("E40").Value = "  +7.28%  "
# Row 41
$ws.Range("D41").Value = "51.30"
$ws.Range("E41").Value = "  +2.30%  "
# Row 42
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").Value = "458.02"
$ws.Range("E42").Value = "  +5.21%  "
# Row 43
$ws.Range("D43").Value = "0.292"
$ws.Range("E43").Value = "  +8.15%  "
# Row 44
$ws.Range("D44").Value = "0.0372"
$ws.Range("E44").Value = "  +1.63%  "
# Row 45
$ws.Range("D45").Value = "2.884.97"
$ws.Range("E45").Value = "  +0.55%  "
# Row 46
$ws.Range("D46").Value = "40.02"
$ws.Range("E46").Value = "  +11.20%  "
# Row 47
$ws.Range("E47").Value = "  -0.23%  "
# Row 48
$ws.Range("D48").Value = "132.98"
$ws.Range("E48").Value = "  +7.47%  "
# Row 50
$ws.Range("E50").Value = "  +0.75%  "
# Row 51
$ws.Range("D51").Value = "2.22"
$ws.Range("E51").Value = "  +3.72%  "

# Restore the default (unstyled) look for column D now that the text is set.
$dRange.Style = "Normal"
